$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.438.34'
$ws.Range("E2").Value = '  +0.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.580.45'
$ws.Range("E3").Value = '  +0.15%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.11'
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.71'
$ws.Range("E8").Value = '  -3.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '23.83'
$ws.Range("E9").Value = '  -0.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.247'
$ws.Range("E10").Value = '  -0.56%  '

$ws.Range("E11").Value = '  -1.08%  '

$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.806.95'
$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.581.01'
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("E15").Value = '  -0.66%  '

$ws.Range("E16").Value = '  -1.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.438.53'
$ws.Range("E17").Value = '  +0.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.95'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.65'
$ws.Range("E19").Value = '  +0.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.44'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0688'
$ws.Range("E21").Value = '  -1.74%  '

$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("E23").Value = '  -0.30%  '

$ws.Range("E24").Value = '  -1.39%  '

$ws.Range("E25").Value = '  +1.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.86'
$ws.Range("E26").Value = '  +0.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.05'
$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.42'
$ws.Range("E28").Value = '  -1.47%  '

$ws.Range("E29").Value = '  -0.84%  '

$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0483'
$ws.Range("E31").Value = '  +3.68%  '

$ws.Range("E32").Value = '  -1.00%  '

$ws.Range("E33").Value = '  -0.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.08'
$ws.Range("E34").Value = '  -1.79%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.400.46'
$ws.Range("E35").Value = '  +0.75%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.10'
$ws.Range("E36").Value = '  +7.75%  '

$ws.Range("E37").Value = '  -3.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  +0.22%  '

$ws.Range("E39").Value = '  +2.16%  '

$ws.Range("E40").Value = '  -1.02%  '

$ws.Range("E41").Value = '  -3.89%  '

$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.89'
$ws.Range("E43").Value = '  +0.81%  '

$ws.Range("E44").Value = '  -1.82%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.42'
$ws.Range("E45").Value = '  -3.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0455'
$ws.Range("E46").Value = '  -2.89%  '

$ws.Range("E47").Value = '  -5.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '62.70'
$ws.Range("E48").Value = '  +0.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.718.78'
$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.07'
$ws.Range("E50").Value = '  +0.08%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  +1.82%  '
